$wb = $excel.ActiveWorkbook
$headerSheet = $wb.Worksheets.Item("header")

# Add the new "Table" sheet right after "header"
$ws = $wb.Worksheets.Add($null, $headerSheet)
$ws.Name = "Table"

# Row 1 - header labels
$ws.Range("A1").Value = "QTY"
$ws.Range("B1").Value = "unit price"
$ws.Range("C1").Value = "total price"

# Row 2 - data, stored as literal text (matches source which keeps
# numeric-looking values as plain strings, not numbers)
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "120,000.00"
$ws.Range("C2").Value = "4120,000.00"
# Drop back to the default (unstyled) look now that the text format has
# done its job of preventing auto-conversion to numbers.
$ws.Range("A2:C2").Style = "Normal"

# Copy the existing bold/centered/bordered header style from "header"!A1
# onto the new header row so it reuses the same style definition.
$headerSheet.Range("A1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null
